# Croatia 1NL.xlsx update script
# Swaps/rotates several data rows (re-sorted by id column B) and
# replaces the final fixture row with a newer record, removing the
# row that is no longer present in the refreshed feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($RowA, $RowB) {
    $rangeA = $ws.Range("B$RowA" + ":AC$RowA")
    $rangeB = $ws.Range("B$RowB" + ":AC$RowB")
    $valA = $rangeA.Value2
    $valB = $rangeB.Value2
    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Simple two-way swaps (everything except column A, the positional index)
Swap-Rows 23 24
Swap-Rows 30 31
Swap-Rows 32 33
Swap-Rows 49 50
Swap-Rows 79 80
Swap-Rows 97 99
Swap-Rows 121 122
Swap-Rows 184 185
Swap-Rows 193 194

# Three-way rotation: new101 = old102, new102 = old103, new103 = old101
Swap-Rows 101 102
Swap-Rows 102 103

# Row 203 receives a brand new fixture record (feed refresh)
$ws.Range("B203").Value2 = 6834829
$ws.Range("E203").Value2 = 45340.36805555555
$ws.Range("F203").Value2 = "Vukovar 91"
$ws.Range("G203").Value2 = "NK Croatia Zmijavci"
$ws.Range("K203").Value2 = 1.571
$ws.Range("L203").Value2 = 3.8
$ws.Range("M203").Value2 = 4.75
$ws.Range("N203").Value2 = 1.571
$ws.Range("O203").Value2 = 3.8
$ws.Range("P203").Value2 = 4.75
$ws.Range("Q203").Value2 = -0.75
$ws.Range("R203").Value2 = 1.75
$ws.Range("S203").Value2 = 2.05
$ws.Range("T203").Value2 = 2.5
$ws.Range("U203").Value2 = 1.9
$ws.Range("V203").Value2 = 1.9
$ws.Range("W203").Value2 = 0
$ws.Range("X203").Value2 = 0
$ws.Range("Y203").Value2 = 0
$ws.Range("Z203").Value2 = 0
$ws.Range("AA203").Value2 = 0

# Row 204 no longer exists in the refreshed feed
$ws.Rows(204).Delete()
